$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 7259.3335
$ws.Cells.Item(43, 9).Value = 5000
$ws.Cells.Item(43, 11).Value = 5000
$ws.Cells.Item(43, 13).Value = -4931

$ws.Cells.Item(74, 8).Value = 2000
$ws.Cells.Item(74, 9).Value = 2000
$ws.Cells.Item(74, 11).Value = 2000
$ws.Cells.Item(74, 13).Value = -1064

$ws.Cells.Item(77, 8).Value = 2000
$ws.Cells.Item(77, 9).Value = 2000
$ws.Cells.Item(77, 11).Value = 10000
$ws.Cells.Item(77, 13).Value = -5320

$ws.Cells.Item(100, 8).Value = 0
$ws.Cells.Item(100, 9).Value = 0
$ws.Cells.Item(100, 11).Value = 0
$ws.Cells.Item(100, 13).Value = ""

$ws.Cells.Item(106, 8).Value = 0
$ws.Cells.Item(106, 9).Value = 0
$ws.Cells.Item(106, 11).Value = 0
$ws.Cells.Item(106, 13).Value = ""

$ws.Cells.Item(113, 8).Value = 17888
$ws.Cells.Item(113, 9).Value = 28776
$ws.Cells.Item(113, 11).Value = 28776
$ws.Cells.Item(113, 13).Value = -25522

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2006.8334
$ws.Cells.Item(2, 10).Value = 3429.4
$ws.Cells.Item(2, 12).Value = 3429.4
$ws.Cells.Item(2, 14).Value = -3655.4

$ws.Cells.Item(32, 8).Value = 3404.6956
$ws.Cells.Item(32, 9).Value = 2882.1904
$ws.Cells.Item(32, 11).Value = 2882.1904
$ws.Cells.Item(32, 13).Value = -2595.1904

$ws.Cells.Item(34, 8).Value = 33999
$ws.Cells.Item(34, 9).Value = 33999
$ws.Cells.Item(34, 11).Value = 33999
$ws.Cells.Item(34, 13).Value = -33728

$ws.Cells.Item(45, 8).Value = 3688.6
$ws.Cells.Item(45, 9).Value = 1443
$ws.Cells.Item(45, 10).Value = 4250
$ws.Cells.Item(45, 11).Value = 1443
$ws.Cells.Item(45, 12).Value = 4250
$ws.Cells.Item(45, 13).Value = -1066
$ws.Cells.Item(45, 14).Value = -5004

$ws.Cells.Item(97, 8).Value = 1916.4
$ws.Cells.Item(97, 9).Value = 1312.25
$ws.Cells.Item(97, 11).Value = 1312.25
$ws.Cells.Item(97, 13).Value = -816.25

$ws.Cells.Item(102, 8).Value = 968
$ws.Cells.Item(102, 9).Value = 968
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 968
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = 654
$ws.Cells.Item(102, 14).Value = ""

$ws.Cells.Item(110, 8).Value = 620.6667
$ws.Cells.Item(110, 9).Value = 620.6667
$ws.Cells.Item(110, 11).Value = 620.6667
$ws.Cells.Item(110, 13).Value = 1424.3333

$ws.Cells.Item(116, 8).Value = 2006.8334
$ws.Cells.Item(116, 10).Value = 3429.4
$ws.Cells.Item(116, 12).Value = 3429.4
$ws.Cells.Item(116, 14).Value = -8017.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2006.8334
$ws.Cells.Item(3, 10).Value = 3429.4
$ws.Cells.Item(3, 12).Value = 3429.4
$ws.Cells.Item(3, 14).Value = -3657.4

$ws.Cells.Item(22, 8).Value = 437.07693
$ws.Cells.Item(22, 9).Value = 437.07693
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 437.07693
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = -264.07693
$ws.Cells.Item(22, 14).Value = ""

$ws.Cells.Item(94, 8).Value = 1222.7778
$ws.Cells.Item(94, 9).Value = 1343.7142
$ws.Cells.Item(94, 10).Value = 799.5
$ws.Cells.Item(94, 11).Value = 1343.7142
$ws.Cells.Item(94, 12).Value = 799.5
$ws.Cells.Item(94, 13).Value = -892.7141999999999
$ws.Cells.Item(94, 14).Value = -1701.5

$ws.Cells.Item(99, 8).Value = 2364.5
$ws.Cells.Item(99, 9).Value = 2038.5
$ws.Cells.Item(99, 11).Value = 2038.5
$ws.Cells.Item(99, 13).Value = -540.5

$ws.Cells.Item(105, 8).Value = 1472
$ws.Cells.Item(105, 9).Value = 1705.25
$ws.Cells.Item(105, 11).Value = 1705.25
$ws.Cells.Item(105, 13).Value = 41.75

$ws.Cells.Item(107, 10).Value = 699
$ws.Cells.Item(107, 12).Value = 699
$ws.Cells.Item(107, 14).Value = -4539

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 3380.5
$ws.Cells.Item(16, 10).Value = 6000
$ws.Cells.Item(16, 12).Value = 6000
$ws.Cells.Item(16, 14).Value = -6574

$ws.Cells.Item(31, 8).Value = 1482.6666
$ws.Cells.Item(31, 10).Value = 1499.5
$ws.Cells.Item(31, 12).Value = 1499.5
$ws.Cells.Item(31, 14).Value = -2089.5

$ws.Cells.Item(34, 8).Value = 1482.6666
$ws.Cells.Item(34, 10).Value = 1499.5
$ws.Cells.Item(34, 12).Value = 1499.5
$ws.Cells.Item(34, 14).Value = -1903.5

$ws.Cells.Item(94, 8).Value = 2410
$ws.Cells.Item(94, 9).Value = 2212
$ws.Cells.Item(94, 11).Value = 2212
$ws.Cells.Item(94, 13).Value = -1761

$ws.Cells.Item(105, 8).Value = 3906.2727
$ws.Cells.Item(105, 9).Value = 1621.125
$ws.Cells.Item(105, 11).Value = 1621.125
$ws.Cells.Item(105, 13).Value = 125.875

$ws.Cells.Item(107, 8).Value = 1407.069
$ws.Cells.Item(107, 9).Value = 1171.0526
$ws.Cells.Item(107, 11).Value = 1171.0526
$ws.Cells.Item(107, 13).Value = 748.9474

$ws.Cells.Item(113, 8).Value = 3380.5
$ws.Cells.Item(113, 10).Value = 6000
$ws.Cells.Item(113, 12).Value = 6000
$ws.Cells.Item(113, 14).Value = -10340

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(75, 8).Value = 1638.5
$ws.Cells.Item(75, 10).Value = 3007.5
$ws.Cells.Item(75, 12).Value = 9022.5
$ws.Cells.Item(75, 14).Value = -11018.5

$ws.Cells.Item(78, 8).Value = 1638.5
$ws.Cells.Item(78, 10).Value = 3007.5
$ws.Cells.Item(78, 12).Value = 27067.5
$ws.Cells.Item(78, 14).Value = -37051.5

$ws.Cells.Item(98, 8).Value = 218.25
$ws.Cells.Item(98, 10).Value = 218
$ws.Cells.Item(98, 12).Value = 654
$ws.Cells.Item(98, 14).Value = -3650

$ws.Cells.Item(112, 8).Value = 2739
$ws.Cells.Item(112, 9).Value = 2913
$ws.Cells.Item(112, 10).Value = 2565
$ws.Cells.Item(112, 11).Value = 8739
$ws.Cells.Item(112, 12).Value = 7695
$ws.Cells.Item(112, 13).Value = -7631
$ws.Cells.Item(112, 14).Value = -9911

$ws.Cells.Item(124, 8).Value = 7205.4
$ws.Cells.Item(124, 9).Value = 3000
$ws.Cells.Item(124, 10).Value = 8256.75
$ws.Cells.Item(124, 11).Value = 9000
$ws.Cells.Item(124, 12).Value = 24770.25
$ws.Cells.Item(124, 13).Value = -4090
$ws.Cells.Item(124, 14).Value = -34590.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(43, 8).Value = 8087.8
$ws.Cells.Item(43, 9).Value = 8756.091
$ws.Cells.Item(43, 10).Value = 6250
$ws.Cells.Item(43, 11).Value = 8756.091
$ws.Cells.Item(43, 12).Value = 6250
$ws.Cells.Item(43, 13).Value = -8605.091
$ws.Cells.Item(43, 14).Value = -6552

$ws.Cells.Item(46, 8).Value = 10241.5
$ws.Cells.Item(46, 9).Value = 5724.5
$ws.Cells.Item(46, 11).Value = 5724.5
$ws.Cells.Item(46, 13).Value = -5568.5

$ws.Cells.Item(57, 8).Value = 5200
$ws.Cells.Item(57, 9).Value = 5200
$ws.Cells.Item(57, 11).Value = 5200
$ws.Cells.Item(57, 13).Value = -4380

$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 13).Value = ""

$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 13).Value = ""

$ws.Cells.Item(97, 8).Value = 1719.8
$ws.Cells.Item(97, 9).Value = 1849.75
$ws.Cells.Item(97, 11).Value = 1849.75
$ws.Cells.Item(97, 13).Value = -1353.75

$ws.Cells.Item(107, 8).Value = 4589.75
$ws.Cells.Item(107, 9).Value = 411.14285
$ws.Cells.Item(107, 11).Value = 411.14285
$ws.Cells.Item(107, 13).Value = 1508.85715

$ws.Cells.Item(132, 8).Value = 3109.6
$ws.Cells.Item(132, 9).Value = 3262
$ws.Cells.Item(132, 11).Value = 9786
$ws.Cells.Item(132, 13).Value = -7256

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2064.6667
$ws.Cells.Item(22, 9).Value = 997.5
$ws.Cells.Item(22, 11).Value = 997.5
$ws.Cells.Item(22, 13).Value = -702.5

$ws.Cells.Item(27, 8).Value = 2064.6667
$ws.Cells.Item(27, 9).Value = 997.5
$ws.Cells.Item(27, 11).Value = 997.5
$ws.Cells.Item(27, 13).Value = -890.5

$ws.Cells.Item(38, 8).Value = 36499.5
$ws.Cells.Item(38, 9).Value = 32999
$ws.Cells.Item(38, 11).Value = 32999
$ws.Cells.Item(38, 13).Value = -32589

$ws.Cells.Item(46, 8).Value = 3693.7805
$ws.Cells.Item(46, 9).Value = 3301.5667
$ws.Cells.Item(46, 11).Value = 3301.5667
$ws.Cells.Item(46, 13).Value = -3113.5667

$ws.Cells.Item(56, 8).Value = 40056
$ws.Cells.Item(56, 9).Value = 40056
$ws.Cells.Item(56, 11).Value = 40056
$ws.Cells.Item(56, 13).Value = -39365

$ws.Cells.Item(63, 8).Value = 41552.5
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 10).Value = 41552.5
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 12).Value = 41552.5
$ws.Cells.Item(63, 13).Value = ""
$ws.Cells.Item(63, 14).Value = -43050.5

$ws.Cells.Item(66, 8).Value = 41552.5
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 10).Value = 41552.5
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 12).Value = 124657.5
$ws.Cells.Item(66, 13).Value = ""
$ws.Cells.Item(66, 14).Value = -132145.5

$ws.Cells.Item(93, 8).Value = 884.4286
$ws.Cells.Item(93, 9).Value = 531.8333
$ws.Cells.Item(93, 11).Value = 531.8333
$ws.Cells.Item(93, 13).Value = 716.1667

$ws.Cells.Item(100, 8).Value = 10586
$ws.Cells.Item(100, 9).Value = 10586
$ws.Cells.Item(100, 11).Value = 10586
$ws.Cells.Item(100, 13).Value = -10045

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(39, 8).Value = 500
$ws.Cells.Item(39, 9).Value = 500
$ws.Cells.Item(39, 11).Value = 500
$ws.Cells.Item(39, 13).Value = -87

$ws.Cells.Item(61, 8).Value = 5025.5
$ws.Cells.Item(61, 9).Value = 5025.5
$ws.Cells.Item(61, 11).Value = 5025.5
$ws.Cells.Item(61, 13).Value = -4733.5

$ws.Cells.Item(96, 8).Value = 1545.8
$ws.Cells.Item(96, 9).Value = 1781.6666
$ws.Cells.Item(96, 10).Value = 1192
$ws.Cells.Item(96, 11).Value = 1781.6666
$ws.Cells.Item(96, 12).Value = 1192
$ws.Cells.Item(96, 13).Value = -408.6666
$ws.Cells.Item(96, 14).Value = -3938
